$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.527.77'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '2.476.64'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.14'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.57'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.544'
$ws.Range("E7").Value = '  -1.22%  '
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.65'
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("E12").Value = '  +2.11%  '
$ws.Range("D13").Value = '2.858.43'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("E14").Value = '  -2.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.18'
$ws.Range("E15").Value = '  +9.14%  '
$ws.Range("D16").Value = '2.446.56'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.766'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").Value = '41.518.52'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.44'
$ws.Range("E19").Value = '  +2.25%  '
$ws.Range("D20").Value = '0.0₃0940'
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.87'
$ws.Range("E21").Value = '  +5.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.25'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.55'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.71'
$ws.Range("E24").Value = '  -1.39%  '
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.80'
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.81'
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.40'
$ws.Range("E31").Value = '  +4.34%  '
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.43'
$ws.Range("E35").Value = '  -8.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.22'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  +4.70%  '
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  -4.01%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.09'
$ws.Range("E41").Value = '  -3.76%  '
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").Value = '1.980.92'
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.26'
$ws.Range("E44").Value = '  -3.71%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.91'
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("D48").Value = '2.719.44'
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.96'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.29'
$ws.Range("E51").Value = '  -2.44%  '
